$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace all occurrences of referee "Sophie Scott" with "Sophie Evelyn"
$used = $ws.UsedRange
$found = $used.Find("Sophie Scott")
if ($found) {
    $firstAddress = $found.Address()
    do {
        $found.Value = "Sophie Evelyn"
        $found = $used.FindNext($found)
    } while ($found -and $found.Address() -ne $firstAddress)
}

# Restore the active selection to A2
$ws.Range("A2").Select()

$wb.Save()
